# Update "想去人数" (want-to-go count) values in column F across sheets
# "展览" (sheet 1), "演出" (sheet 2) and "全部类型" (sheet 4).

$wb = $excel.ActiveWorkbook

# --- 展览 sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 274
$ws1.Range("F5").Value  = 2877
$ws1.Range("F8").Value  = 2233
$ws1.Range("F9").Value  = 1392
$ws1.Range("F11").Value = 441
$ws1.Range("F16").Value = 4773
$ws1.Range("F18").Value = 5256
$ws1.Range("F19").Value = 1811
$ws1.Range("F29").Value = 2002
$ws1.Range("F34").Value = 357

# --- 演出 sheet ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value  = 102
$ws2.Range("F18").Value = 50

# --- 全部类型 sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 102
$ws4.Range("F9").Value  = 274
$ws4.Range("F10").Value = 2877
$ws4.Range("F12").Value = 2233
$ws4.Range("F13").Value = 1392
$ws4.Range("F17").Value = 441
$ws4.Range("F25").Value = 4773
$ws4.Range("F27").Value = 5256
$ws4.Range("F28").Value = 1811
$ws4.Range("F42").Value = 50
$ws4.Range("F43").Value = 2002
$ws4.Range("F48").Value = 357
